$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "dal + roshun = 24" -> increase Bazar TK for day 19 (column T, row 43) by 24 (10 -> 34)
$ws.Range("T43").Value = 34

# Update the active selection to reflect where the author ended up (T44)
$ws.Range("T44").Select()
